$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'37.409.74"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.92%  '
$ws.Range('D3').Value = "'2.065.06"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = "'231.59"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('D6').Value = "'0.627"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.65%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = "'57.01"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.55%  '
$ws.Range('D9').Value = "'0.387"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.63%  '
$ws.Range('D10').Value = "'0.0775"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('D11').Value = "'0.108"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.75%  '
$ws.Range('D12').Value = "'14.82"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('D13').Value = "'2.374.71"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.23%  '
$ws.Range('D14').Value = "'20.81"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').Value = "'0.762"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.47%  '
$ws.Range('D16').Value = "'5.30"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('D17').Value = "'2.075.41"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.59%  '
$ws.Range('D18').Value = "'37.347.02"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').Value = "'70.28"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.11%  '
$ws.Range('D20').Value = "'5.95"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.90%  '
$ws.Range('D21').Value = "'0.0₃0826"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.96%  '
$ws.Range('D22').Value = "'227.68"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('E24').Value = '  -1.29%  '
$ws.Range('D25').Value = "'2.36"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.27%  '
$ws.Range('D26').Value = "'9.58"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.22%  '
$ws.Range('D27').Value = "'169.89"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('E28').Value = '  -4.37%  '
$ws.Range('D29').Value = "'19.41"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('E30').Value = '  -1.60%  '
$ws.Range('D31').Value = "'0.122"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.85%  '
$ws.Range('D32').Value = "'4.60"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.80%  '
$ws.Range('D33').Value = "'0.0629"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.55%  '
$ws.Range('D34').Value = "'4.61"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.90%  '
$ws.Range('D35').Value = "'2.45"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.24%  '
$ws.Range('D36').Value = "'1.82"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('D37').Value = "'3.29"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.46%  '
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('D39').Value = "'5.26"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.97%  '
$ws.Range('D40').Value = "'0.0229"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.56%  '
$ws.Range('D41').Value = "'99.37"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.04%  '
$ws.Range('E42').Value = '  +0.96%  '
$ws.Range('D43').Value = "'0.0950"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.48%  '
$ws.Range('D44').Value = "'1.19"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.77%  '
$ws.Range('D45').Value = "'1.458.14"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.90%  '
$ws.Range('D46').Value = "'16.61"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.40%  '
$ws.Range('E47').Value = '  -2.11%  '
$ws.Range('D48').Value = "'3.94"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.05%  '
$ws.Range('E49').Value = '  -2.72%  '
$ws.Range('D50').Value = "'2.93"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.40%  '
$ws.Range('D51').Value = "'2.260.96"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.20%  '
